# Auto-generated: apply updated crypto price/volume values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.133.94"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = "'1.900.42"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'307.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').Value = "'0.5233"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('D9').Value = "'0.07287"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').Value = "'21.38"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = "'0.08204"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').Value = "'95.53"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = "'1.846.16"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('D15').Value = "'5.355"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = "'0.000008664"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = "'14.69"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = "'27.177.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = "'2.098.07"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('D23').Value = "'10.79"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('D24').Value = "'6.451"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = "'149.76"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.26%  '
$ws.Range('D26').Value = "'2.320"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = "'1.740"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').Value = "'115.66"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').Value = "'4.901"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = "'0.09220"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').Value = "'0.05042"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = "'0.7937"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('D35').Value = "'1.222"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('D36').Value = "'2.961"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('D37').Value = "'3.360"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').Value = "'2.628"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('D39').Value = "'0.5730"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').Value = "'0.01991"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').Value = "'1.082"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').Value = "'9.118"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('D43').Value = "'6.621"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D44').Value = "'116.37"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').Value = "'0.1517"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = "'0.4896"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = "'10.14"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('D50').Value = "'38.52"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('D51').Value = "'64.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.76%  '
